# Delete the post row containing "「自分と仲直りしよう、一杯の珈琲で！」" (row 490).
# All subsequent rows shift up by one, and the sheet's used range shrinks by one row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(490).Delete()
